$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update metric values
$ws.Range("B3").Value = 10685437837766.5
$ws.Range("C3").Value = 8679957574858.971
$ws.Range("D3").Value = 8882596980697.797

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 6499717772070.138
$ws.Range("C4").Value = 6499717772070.137
$ws.Range("D4").Value = 6656337477421.225

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 38858527306381.99
$ws.Range("C5").Value = 56358530489483.09
$ws.Range("D5").Value = 67463987872022.22
